$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function CopyFmt($srcAddr, $dstAddr) {
  $ws.Range($srcAddr).Copy()
  $ws.Range($dstAddr).PasteSpecial(-4122)
  $excel.CutCopyMode = $false
}

# ---- Row 43 : copy cell formatting from existing rows that already carry ----
# ---- the exact style needed, then overwrite with the new literal values  ----
CopyFmt "A41" "A43"
CopyFmt "B40" "B43"
CopyFmt "C41" "C43"
CopyFmt "D41" "D43"
CopyFmt "E41" "E43"
CopyFmt "F40" "F43"
CopyFmt "G41" "G43"
CopyFmt "H41" "H43"
CopyFmt "I41" "I43"

$ws.Range("A43").Value = 44691
$ws.Range("B43").Value = "BUAM"
$ws.Range("C43").Value = "N/A"
$ws.Range("D43").Value = "Saint-Ignace-de-Loyola"
$ws.Range("E43").Value = "Montérégie/Lanaudière"
$ws.Range("F43").Value = "A"
$ws.Range("G43").Value = "Cote 1"
$ws.Range("H43").Value = "iNaturalist"
$ws.Range("I43").Value = "Timothe Breton"

# ---- Row 44 ----
CopyFmt "A41" "A44"
CopyFmt "B40" "B44"
CopyFmt "C41" "C44"
CopyFmt "D41" "D44"
CopyFmt "E41" "E44"
CopyFmt "F11" "F44"
CopyFmt "G41" "G44"
CopyFmt "H41" "H44"
CopyFmt "I41" "I44"

$ws.Range("A44").Value = 44691
$ws.Range("B44").Value = "BUAM"
$ws.Range("C44").Value = "N/A"
$ws.Range("D44").Value = "La Conception"
$ws.Range("E44").Value = "Laurentides"
$ws.Range("F44").Value = "B"
$ws.Range("G44").Value = "Cote 2"
$ws.Range("H44").Value = "iNaturalist"
$ws.Range("I44").Value = "Caro Marcotte"

# ---- Selection, matching the recorded cursor position after the edit ----
$ws.Range("H49").Select()
